$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.519.27"
$ws.Range("E2").Value = "  +0.03%  "

$ws.Range("D3").Value = "3.062.90"
$ws.Range("E3").Value = "  +2.62%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "386.05"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.02%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "102.95"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.69%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.543"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.49%  "

$ws.Range("E9").Value = "  -1.34%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.76"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.00%  "

$ws.Range("E11").Value = "  +0.17%  "

$ws.Range("E12").Value = "  +0.30%  "

$ws.Range("D13").Value = "3.554.15"
$ws.Range("E13").Value = "  +3.15%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.63"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.97%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.77"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.23%  "

$ws.Range("D16").Value = "3.070.20"
$ws.Range("E16").Value = "  +2.47%  "

$ws.Range("E17").Value = "  -2.15%  "

$ws.Range("E18").Value = "  -3.75%  "

$ws.Range("D19").Value = "51.621.36"
$ws.Range("E19").Value = "  +0.25%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.15"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.88%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.43"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.59%  "

$ws.Range("D22").Value = "0.0₃0968"
$ws.Range("E22").Value = "  +0.44%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.17"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.43%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "267.77"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.11%  "

$ws.Range("E25").Value = "  -1.90%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.16"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.74%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "26.89"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.06%  "

$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.28"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.77%  "

$ws.Range("B29").Value = "Kaspa"
$ws.Range("C29").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.170"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.94%  "

$ws.Range("E30").Value = "  +0.05%  "

$ws.Range("E31").Value = "  -2.23%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "10.25"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.81%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "34.72"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.20%  "

$ws.Range("E34").Value = "  -0.01%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "50.11"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.82%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0448"
$ws.Range("D36").Style = "Normal"

$ws.Range("E37").Value = "  -0.02%  "

$ws.Range("E38").Value = "  +1.73%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.293"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +8.29%  "

$ws.Range("E40").Value = "  +1.50%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "16.88"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.42%  "

$ws.Range("E42").Value = "  +0.55%  "

$ws.Range("E43").Value = "  -0.90%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "125.30"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.25%  "

$ws.Range("E45").Value = "  -0.17%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "21.92"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.24%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.09"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.61%  "

$ws.Range("E48").Value = "  +2.53%  "

$ws.Range("D49").Value = "2.034.42"
$ws.Range("E49").Value = "  +0.03%  "

$ws.Range("D50").Value = "3.366.75"
$ws.Range("E50").Value = "  +2.62%  "

$ws.Range("E51").Value = "  +7.22%  "
